$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear "KW 32" values from B3 and B4, and set it on B5 instead
$ws.Range("B3").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("B5").Value = "KW 32"

# Rearrange column A values/styles for rows 8-11
$ws.Range("A8").Value = "Interessenprofileverwaltung"
$ws.Range("A9").Clear()
$ws.Range("A10").Value = "Soll Kriterien:"
$ws.Range("A10").Font.Bold = $true
$ws.Range("A11").Value = "Löschen von Interessenprofilen"

# Update the selected cell shown when the workbook is opened
$ws.Range("B8").Select()
